# LDAP Client Linux finish
# Record the completion of the "1 Windows-Client & 1 Linux-Client" task:
# the Linux client sub-row (row 30) gets the responsible person, the
# completion date and the time spent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wer(Adler, Karic, Kopec) -> Adler did the work
$ws.Range("C30").Value = "Adler"

# Datum -> 2015-03-06 (same date style as the other "Datum" cells, e.g. D18)
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("D30").Value = 42069

# Zeitaufwand -> 180min
$ws.Range("E30").Value = "180min"

# Leave the cursor where the author left it when saving
$ws.Range("D42").Select() | Out-Null
